$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; remaining columns (B:F) shift left to (A:E),
# matching the target layout.
$ws.Range("A1").EntireColumn.Delete()
